$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1271.1428
$ws.Range("I28").Value = 1392.1818
$ws.Range("K28").Value = 1392.1818
$ws.Range("M28").Value = -907.1818000000001
$ws.Range("H80").Value = 519.9643
$ws.Range("I80").Value = 467.81818
$ws.Range("J80").Value = 553.7059
$ws.Range("K80").Value = 1403.45454
$ws.Range("L80").Value = 1661.1177
$ws.Range("M80").Value = -405.45454
$ws.Range("N80").Value = -3657.1177
$ws.Range("H83").Value = 519.9643
$ws.Range("I83").Value = 467.81818
$ws.Range("J83").Value = 553.7059
$ws.Range("K83").Value = 4210.36362
$ws.Range("L83").Value = 4983.3531
$ws.Range("M83").Value = 781.6363799999999
$ws.Range("N83").Value = -14967.3531
$ws.Range("H88").Value = 861.6875
$ws.Range("I88").Value = 775.375
$ws.Range("J88").Value = 948
$ws.Range("K88").Value = 775.375
$ws.Range("L88").Value = 948
$ws.Range("M88").Value = -369.375
$ws.Range("N88").Value = -1760
$ws.Range("H91").Value = 861.6875
$ws.Range("I91").Value = 775.375
$ws.Range("J91").Value = 948
$ws.Range("K91").Value = 775.375
$ws.Range("L91").Value = 948
$ws.Range("M91").Value = 628.625
$ws.Range("N91").Value = -3756
$ws.Range("H107").Value = 5796.7144
$ws.Range("I107").Value = 1068.1666
$ws.Range("K107").Value = 1068.1666
$ws.Range("M107").Value = 851.8334
$ws.Range("H111").Value = 71429576
$ws.Range("I111").Value = 83334480
$ws.Range("J111").Value = 150
$ws.Range("K111").Value = 250003440
$ws.Range("L111").Value = 450
$ws.Range("M111").Value = -250000373
$ws.Range("N111").Value = -6584
$ws.Range("H112").Value = 2410
$ws.Range("J112").Value = 2962.5
$ws.Range("L112").Value = 8887.5
$ws.Range("N112").Value = -11103.5
$ws.Range("H116").Value = 2851910.8
$ws.Range("I116").Value = 6995672
$ws.Range("J116").Value = 3075
$ws.Range("K116").Value = 6995672
$ws.Range("L116").Value = 3075
$ws.Range("M116").Value = -6992230
$ws.Range("N116").Value = -9959
$ws.Range("H129").Value = 996.62335
$ws.Range("I129").Value = 350
$ws.Range("J129").Value = 1032.0548
$ws.Range("K129").Value = 1050
$ws.Range("L129").Value = 3096.1644
$ws.Range("M129").Value = 3950
$ws.Range("N129").Value = -13096.1644

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2384.5
$ws.Range("I2").Value = 2470.2222
$ws.Range("K2").Value = 2470.2222
$ws.Range("M2").Value = -2357.2222
$ws.Range("H45").Value = 4093.6
$ws.Range("I45").Value = 4012
$ws.Range("K45").Value = 4012
$ws.Range("M45").Value = -3635
$ws.Range("H110").Value = 1827.5555
$ws.Range("I110").Value = 1149.6
$ws.Range("J110").Value = 2675
$ws.Range("K110").Value = 1149.6
$ws.Range("L110").Value = 2675
$ws.Range("M110").Value = 895.4000000000001
$ws.Range("N110").Value = -6765
$ws.Range("H116").Value = 2384.5
$ws.Range("I116").Value = 2470.2222
$ws.Range("K116").Value = 2470.2222
$ws.Range("M116").Value = -176.2222000000002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2384.5
$ws.Range("I3").Value = 2470.2222
$ws.Range("K3").Value = 2470.2222
$ws.Range("M3").Value = -2356.2222
$ws.Range("H86").Value = 2704.5454
$ws.Range("I86").Value = 2125
$ws.Range("K86").Value = 2125
$ws.Range("M86").Value = -1002
$ws.Range("H89").Value = 2704.5454
$ws.Range("I89").Value = 2125
$ws.Range("K89").Value = 10625
$ws.Range("M89").Value = -5009
$ws.Range("H107").Value = 2831.8572
$ws.Range("I107").Value = 3001.6667
$ws.Range("K107").Value = 3001.6667
$ws.Range("M107").Value = -1081.6667

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2757
$ws.Range("I16").Value = 1883.1666
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 1883.1666
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -1596.1666
$ws.Range("N16").Value = -8574
$ws.Range("H107").Value = 746.5789
$ws.Range("I107").Value = 530
$ws.Range("J107").Value = 1117.8572
$ws.Range("K107").Value = 530
$ws.Range("L107").Value = 1117.8572
$ws.Range("M107").Value = 1390
$ws.Range("N107").Value = -4957.8572
$ws.Range("H113").Value = 2757
$ws.Range("I113").Value = 1883.1666
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 1883.1666
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = 286.8334
$ws.Range("N113").Value = -12340

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 13.3125
$ws.Range("I12").Value = 3.75
$ws.Range("J12").Value = 16.5
$ws.Range("K12").Value = 11.25
$ws.Range("L12").Value = 49.5
$ws.Range("M12").Value = 161.75
$ws.Range("N12").Value = -395.5
$ws.Range("H131").Value = 2481.082
$ws.Range("I131").Value = 539.0909
$ws.Range("J131").Value = 2908.32
$ws.Range("K131").Value = 1617.2727
$ws.Range("L131").Value = 8724.960000000001
$ws.Range("M131").Value = 3422.7273
$ws.Range("N131").Value = -18804.96

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9998
$ws.Range("I70").Value = 11147.5
$ws.Range("J70").Value = 5400
$ws.Range("K70").Value = 11147.5
$ws.Range("L70").Value = 5400
$ws.Range("M70").Value = -10877.5
$ws.Range("N70").Value = -5940
$ws.Range("H73").Value = 9998
$ws.Range("I73").Value = 11147.5
$ws.Range("J73").Value = 5400
$ws.Range("K73").Value = 11147.5
$ws.Range("L73").Value = 5400
$ws.Range("M73").Value = -10211.5
$ws.Range("N73").Value = -7272
$ws.Range("H107").Value = 206.90909
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H113").Value = 19866.5
$ws.Range("I113").Value = 2799.75
$ws.Range("K113").Value = 2799.75
$ws.Range("M113").Value = -629.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 938.875
$ws.Range("I16").Value = 858.7143
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 858.7143
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -688.7143
$ws.Range("N16").Value = -1840
$ws.Range("H55").Value = 257.12903
$ws.Range("I55").Value = 207.65
$ws.Range("J55").Value = 347.0909
$ws.Range("K55").Value = 207.65
$ws.Range("L55").Value = 347.0909
$ws.Range("M55").Value = -34.65000000000001
$ws.Range("N55").Value = -693.0908999999999
$ws.Range("H61").Value = 20698.6
$ws.Range("I61").Value = 33664.668
$ws.Range("J61").Value = 1249.5
$ws.Range("K61").Value = 33664.668
$ws.Range("L61").Value = 1249.5
$ws.Range("M61").Value = -33462.668
$ws.Range("N61").Value = -1653.5
$ws.Range("H113").Value = 20698.6
$ws.Range("I113").Value = 33664.668
$ws.Range("J113").Value = 1249.5
$ws.Range("K113").Value = 33664.668
$ws.Range("L113").Value = 1249.5
$ws.Range("M113").Value = -31494.668
$ws.Range("N113").Value = -5589.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 357
$ws.Range("I107").Value = 332.7
$ws.Range("K107").Value = 998.0999999999999
$ws.Range("M107").Value = 921.9000000000001
$ws.Range("H113").Value = 273.1111
$ws.Range("I113").Value = 261.5625
$ws.Range("K113").Value = 784.6875
$ws.Range("M113").Value = 1385.3125
